$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 5 (2013年) values
$ws.Range("B5").Value = 790297
$ws.Range("C5").Value = 8517893
$ws.Range("D5").Value = 1437073
$ws.Range("F5").Value = 262434

# Add new row 12 (2021年): copy formatting from row 11's A cell, then fill values
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 741801
$ws.Range("C12").Value = 28665212
$ws.Range("D12").Value = 2900264
$ws.Range("E12").Value = 223498
$ws.Range("F12").Value = 336197
